$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 105 (the existing data from row 105 onward shifts down to 109+)
$ws.Rows("105:108").Insert()

# Row 105
$ws.Range("A105").Value = 6
$ws.Range("B105").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C105").Value = "Metropolitana"
$ws.Range("D105").Value = 44476
$ws.Range("E105").Value = 13
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100107
$ws.Range("H105").Value = "Otros"
$ws.Range("I105").Value = 100107002
$ws.Range("J105").Value = "Chirimoya"
$ws.Range("K105").Value = "Cultivar IV Región"
$ws.Range("L105").Value = "Especial"
$ws.Range("M105").Value = 100
$ws.Range("N105").Value = 2700
$ws.Range("O105").Value = 2700
$ws.Range("P105").Value = 2700
$ws.Range("Q105").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R105").Value = "Provincia de Limarí"
$ws.Range("S105").Value = 2700
$ws.Range("T105").Value = 1

# Row 106
$ws.Range("A106").Value = 6
$ws.Range("B106").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C106").Value = "Metropolitana"
$ws.Range("D106").Value = 44476
$ws.Range("E106").Value = 13
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100107
$ws.Range("H106").Value = "Otros"
$ws.Range("I106").Value = 100107002
$ws.Range("J106").Value = "Chirimoya"
$ws.Range("K106").Value = "Cultivar IV Región"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 150
$ws.Range("N106").Value = 2500
$ws.Range("O106").Value = 2500
$ws.Range("P106").Value = 2500
$ws.Range("Q106").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R106").Value = "Provincia de Limarí"
$ws.Range("S106").Value = 2500
$ws.Range("T106").Value = 1

# Row 107
$ws.Range("A107").Value = 6
$ws.Range("B107").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C107").Value = "Metropolitana"
$ws.Range("D107").Value = 44476
$ws.Range("E107").Value = 13
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100107
$ws.Range("H107").Value = "Otros"
$ws.Range("I107").Value = 100107002
$ws.Range("J107").Value = "Chirimoya"
$ws.Range("K107").Value = "Cultivar IV Región"
$ws.Range("L107").Value = "Segunda"
$ws.Range("M107").Value = 150
$ws.Range("N107").Value = 2200
$ws.Range("O107").Value = 2200
$ws.Range("P107").Value = 2200
$ws.Range("Q107").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R107").Value = "Provincia de Limarí"
$ws.Range("S107").Value = 2200
$ws.Range("T107").Value = 1

# Row 108
$ws.Range("A108").Value = 6
$ws.Range("B108").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C108").Value = "Metropolitana"
$ws.Range("D108").Value = 44476
$ws.Range("E108").Value = 13
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100107
$ws.Range("H108").Value = "Otros"
$ws.Range("I108").Value = 100107002
$ws.Range("J108").Value = "Chirimoya"
$ws.Range("K108").Value = "Cultivar IV Región"
$ws.Range("L108").Value = "Tercera"
$ws.Range("M108").Value = 100
$ws.Range("N108").Value = 1800
$ws.Range("O108").Value = 1800
$ws.Range("P108").Value = 1800
$ws.Range("Q108").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R108").Value = "Provincia de Limarí"
$ws.Range("S108").Value = 1800
$ws.Range("T108").Value = 1
